# Update "Ngày công" (work days) and the dependent salary figures on the
# "Lương" sheet: công tăng từ 18 lên 19 ngày, kéo theo phụ cấp và lương.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 19
$ws.Range("B3").Value = 665000
$ws.Range("B12").Value = 2035714.285714286
$ws.Range("B29").Value = 2900714.285714285
$ws.Range("B31").Value = 2900714.285714285
